# Updates cryptos list: refreshed Price (D) and Volume(1h) (E) columns,
# plus a Coin/Link swap between rows 40-41 (ImmutableX <-> MXToken).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.680.87"
$ws.Range("E2").Value = "  +1.35%  "

$ws.Range("D3").Value = "1.634.11"
$ws.Range("E3").Value = "  +1.25%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("E7").Value = "  +1.01%  "

$ws.Range("D8").Value = "0.251"
$ws.Range("E8").Value = "  +0.29%  "

$ws.Range("E9").Value = "  +0.64%  "

$ws.Range("D10").Value = "19.01"
$ws.Range("E10").Value = "  +2.77%  "

$ws.Range("D11").Value = "0.0834"
$ws.Range("E11").Value = "  +2.45%  "

$ws.Range("D12").Value = "1.864.45"
$ws.Range("E12").Value = "  +1.68%  "

$ws.Range("D13").Value = "1.633.20"
$ws.Range("E13").Value = "  +2.28%  "

$ws.Range("D14").Value = "4.05"
$ws.Range("E14").Value = "  -0.11%  "

$ws.Range("D15").Value = "0.524"
$ws.Range("E15").Value = "  +1.38%  "

$ws.Range("D16").Value = "26.663.48"
$ws.Range("E16").Value = "  +1.33%  "

$ws.Range("D17").Value = "63.12"
$ws.Range("E17").Value = "  +1.87%  "

$ws.Range("D18").Value = "0.0₃0733"

$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.05%  "

$ws.Range("D21").Value = "4.33"
$ws.Range("E21").Value = "  +0.39%  "

$ws.Range("D22").Value = "9.38"
$ws.Range("E22").Value = "  +0.60%  "

$ws.Range("D23").Value = "6.08"
$ws.Range("E23").Value = "  +0.97%  "

$ws.Range("E24").Value = "  -1.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.03%  "

$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("E27").Value = "  -1.88%  "

$ws.Range("D28").Value = "15.36"
$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("D29").Value = "6.67"
$ws.Range("E29").Value = "  +1.63%  "

$ws.Range("E30").Value = "  +5.47%  "

$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("D33").Value = "2.94"
$ws.Range("E33").Value = "  -0.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.92%  "

$ws.Range("E35").Value = "  -0.50%  "

$ws.Range("D36").Value = "1.163.91"
$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("E37").Value = "  +0.63%  "

$ws.Range("D38").Value = "0.811"
$ws.Range("E38").Value = "  +1.59%  "

$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.32"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").Value = "0.502"
$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.71%  "

$ws.Range("D43").Value = "0.793"
$ws.Range("E43").Value = "  +0.83%  "

$ws.Range("D44").Value = "1.773.75"
$ws.Range("E44").Value = "  +1.50%  "

$ws.Range("D45").Value = "92.42"
$ws.Range("E45").Value = "  +0.60%  "

$ws.Range("D46").Value = "1.54"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").Value = "54.65"
$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("E48").Value = "  +0.75%  "

$ws.Range("D49").Value = "7.59"
$ws.Range("E49").Value = "  +4.09%  "

$ws.Range("D50").Value = "0.409"
$ws.Range("E50").Value = "  +0.69%  "

$ws.Range("E51").Value = "  +0.04%  "
